# Adds two new columns ("BL" and "Operating Freq") to the circuit-component
# table on Hoja1, with a value of 0 for the existing data row, matching the
# style (centered alignment) used by the other header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1, G1) -- text values, centered like the rest of row 1.
$ws.Range("F1").Value = "BL"
$ws.Range("G1").Value = "Operating Freq"

$headerRange = $ws.Range("F1:G1")
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4108     # xlCenter

# New data cells (F2, G2) for the existing row.
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

# Match the author's last-saved selection (cell F2).
[void]$ws.Range("F2").Select()
